$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 10379.42
$ws.Range("B8").Value = 10488.5
$ws.Range("C8").Value = 109.08
$ws.Range("D8").Value = 107.95
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = -1.04
$ws.Range("G8").Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(42612.672962962963)
$ws.Range("H8").Value = $false
